$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Paragraph 2: drop "Apparently, this is compressing… (it's not)."
#    leaving "...words positions being stored in another file."
#    then re-split the remaining sentence into two runs at
#    "position" | "s being stored in another file."
# ------------------------------------------------------------------
$d.Content.Find.Execute(" Apparently, this is compressing… (it’s not).", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$d.Content.Find.Execute("positions being stored", $true, $false, $false, $false, $false, $true, 1, $false, "position", 2)

$findRng = $d.Content
$findRng.Find.Execute("position", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $findRng.End
$splitRng = $d.Range($splitPos, $splitPos)
$splitRng.InsertAfter("s being stored")

# Force a real run boundary at $splitPos (a transient bookmark leaves
# the run split behind even once it is removed again).
$tempBmRng = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("TEMP_SPLIT_1", $tempBmRng)
$d.Bookmarks("TEMP_SPLIT_1").Delete()

# ------------------------------------------------------------------
# 2. Paragraph starting "Let's take a look..." -> split into
#    "Let's" | " take " | (bookmark _GoBack) | "a look at ..."
# ------------------------------------------------------------------
$letsRng = $d.Content
$letsRng.Find.Execute("Let’s take a look at the criteria list from the task analysis for task 2.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$letsStart = $letsRng.Start

$splitA = $letsStart + 5     # after "Let's"
$splitB = $letsStart + 11    # after "Let's take "

$rngA = $d.Range($splitA, $splitA)
$d.Bookmarks.Add("TEMP_SPLIT_2", $rngA)
$d.Bookmarks("TEMP_SPLIT_2").Delete()

$rngB = $d.Range($splitB, $splitB)
$d.Bookmarks.Add("TEMP_SPLIT_3", $rngB)
$d.Bookmarks("TEMP_SPLIT_3").Delete()

# ------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the end of the document
#    (after "capitalisation and punctuation") to between
#    "Let's take " and "a look at ...".
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$goBackRng = $d.Range($splitB, $splitB)
$d.Bookmarks.Add("_GoBack", $goBackRng)
